# Applies the cryptos-list price/volume refresh described in the commit diff.
# Numeric-looking text in column D must stay TEXT (source data is inlineStr),
# so those assignments use a leading apostrophe (quote-prefix), matching how
# Excel keeps user-typed numeric strings as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.523.03"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "2.467.21"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'313.68"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").Value = "'91.61"
$ws.Range("E6").Value = "  -2.74%  "

$ws.Range("E7").Value = "  -0.45%  "

$ws.Range("E8").Value = "  -0.17%  "

$ws.Range("E9").Value = "  +2.68%  "

$ws.Range("D10").Value = "'32.29"
$ws.Range("E10").Value = "  -3.46%  "

$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  +1.02%  "

$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("D13").Value = "2.849.34"
$ws.Range("E13").Value = "  -0.27%  "

$ws.Range("E14").Value = "  -2.10%  "

$ws.Range("D15").Value = "'16.05"
$ws.Range("E15").Value = "  +5.05%  "

$ws.Range("D16").Value = "2.503.02"
$ws.Range("E16").Value = "  +1.38%  "

$ws.Range("D17").Value = "'0.767"
$ws.Range("E17").Value = "  -2.77%  "

$ws.Range("D18").Value = "41.493.02"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("E19").Value = "  +2.97%  "

$ws.Range("D20").Value = "0.0₃0945"
$ws.Range("E20").Value = "  +2.17%  "

$ws.Range("D21").Value = "'71.51"
$ws.Range("E21").Value = "  +4.12%  "

$ws.Range("D22").Value = "'11.08"
$ws.Range("E22").Value = "  -1.17%  "

$ws.Range("D23").Value = "'236.06"
$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  -0.97%  "

$ws.Range("D27").Value = "'24.70"
$ws.Range("E27").Value = "  +2.80%  "

$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "'35.34"
$ws.Range("E30").Value = "  -2.52%  "

$ws.Range("D31").Value = "'156.33"
$ws.Range("E31").Value = "  +3.11%  "

$ws.Range("D32").Value = "'5.43"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("E33").Value = "  +0.28%  "

$ws.Range("D34").Value = "'0.0756"
$ws.Range("E34").Value = "  +1.52%  "

$ws.Range("D35").Value = "'17.17"
$ws.Range("E35").Value = "  -0.71%  "

$ws.Range("D36").Value = "'2.34"
$ws.Range("E36").Value = "  -8.81%  "

$ws.Range("D37").Value = "'2.87"
$ws.Range("E37").Value = "  -6.60%  "

$ws.Range("E38").Value = "  +1.80%  "

$ws.Range("E39").Value = "  -0.67%  "

$ws.Range("E40").Value = "  -4.56%  "

$ws.Range("D41").Value = "'4.03"
$ws.Range("E41").Value = "  -3.89%  "

$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").Value = "1.952.34"
$ws.Range("E43").Value = "  -1.83%  "

$ws.Range("D44").Value = "'0.0283"
$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "'18.70"
$ws.Range("E45").Value = "  -3.27%  "

$ws.Range("E46").Value = "  -2.33%  "

$ws.Range("D47").Value = "'9.06"
$ws.Range("E47").Value = "  +3.93%  "

$ws.Range("D48").Value = "2.709.50"
$ws.Range("E48").Value = "  -0.41%  "

$ws.Range("D49").Value = "'96.96"
$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").Value = "'66.80"
$ws.Range("E50").Value = "  -3.64%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.171"
$ws.Range("E51").Value = "  -3.50%  "
